$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.232819676399231
$ws.Range("B1").Value = 1.895875096321106
$ws.Range("C1").Value = 4.19976806640625
$ws.Range("D1").Value = 3.165406942367554
$ws.Range("E1").Value = 1.178413510322571
